# Generate Report for Handoff
# Updates status/date strings to reflect a fresh handoff, and resizes
# the "Status" / status-like columns on each sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-07 07:24:09"

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-07 07:23:57"

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-07 07:24:09"

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
